# Apply the commit's changes to the Dombrovskis meetings sheet.
#
# Net effect of the diff is a swap of data between three pairs of rows
# (the date in column A is identical within each pair, so only the other
# columns actually move):
#   - Row 44  <-> Row 45  (columns C, D)
#   - Row 49  <-> Row 50  (columns B, C, D)
#   - Row 52  <-> Row 53  (columns B, C, D)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Range {
    param($ws, [string]$addr1, [string]$addr2)
    $r1 = $ws.Range($addr1)
    $r2 = $ws.Range($addr2)
    $v1 = $r1.Value()
    $v2 = $r2.Value()
    $r1.Value = $v2
    $r2.Value = $v1
}

# --- Rows 44 / 45: swap Entity (C) and Subject (D) ---
Swap-Range $ws "C44" "C45"
Swap-Range $ws "D44" "D45"

# --- Rows 49 / 50: swap Location (B), Entity (C) and Subject (D) ---
Swap-Range $ws "B49" "B50"
Swap-Range $ws "C49" "C50"
Swap-Range $ws "D49" "D50"

# --- Rows 52 / 53: swap Location (B), Entity (C) and Subject (D) ---
Swap-Range $ws "B52" "B53"
Swap-Range $ws "C52" "C53"
Swap-Range $ws "D52" "D53"

# Writing new text into the wrapped cells above makes Excel recompute the
# row height as if it were a custom height; put it back the way it was
# (rows had no explicit custom height in the original workbook).
$ws.Rows.Item(44).EntireRow.AutoFit()
$ws.Rows.Item(45).EntireRow.AutoFit()
$ws.Rows.Item(49).EntireRow.AutoFit()
$ws.Rows.Item(50).EntireRow.AutoFit()
$ws.Rows.Item(52).EntireRow.AutoFit()
$ws.Rows.Item(53).EntireRow.AutoFit()
